# Apply the new stimuli order: shuffled image/word/category shared-string
# content and a new per-row count (column B), matching the commit
# "initial version of stimuli order".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: count (B), image (C), word (D), category (E)
$rows = @(
    @(66, "dog/dog107.png", "wecken", "dog"),
    @(25, "car/car095.png", "zögern", "car"),
    @(6, "car/car080.png", "atmen", "car"),
    @(103, "dog/dog118.png", "rufen", "dog"),
    @(50, "dog/dog088.png", "frischen", "dog"),
    @(105, "dog/dog089.png", "mögen", "dog"),
    @(20, "dog/dog071.png", "lernen", "dog"),
    @(98, "dog/dog077.png", "hören", "dog"),
    @(65, "dog/dog081.png", "wehen", "dog"),
    @(10, "car/car074.png", "spüren", "car"),
    @(16, "car/car101.png", "binden", "car"),
    @(39, "car/car064.png", "rechnen", "car"),
    @(34, "car/car093.png", "narren", "car"),
    @(14, "dog/dog072.png", "danken", "dog"),
    @(84, "car/car075.png", "süßen", "car"),
    @(97, "car/car077.png", "spenden", "car"),
    @(60, "car/car111.png", "bergen", "car"),
    @(12, "dog/dog085.png", "proben", "dog"),
    @(37, "dog/dog095.png", "lächeln", "dog"),
    @(90, "dog/dog076.png", "leeren", "dog"),
    @(49, "dog/dog100.png", "legen", "dog"),
    @(15, "car/car086.png", "duschen", "car"),
    @(86, "car/car115.png", "stoppen", "car"),
    @(36, "car/car076.png", "hacken", "car"),
    @(127, "car/car105.png", "kriegen", "car"),
    @(106, "dog/dog065.png", "meinen", "dog"),
    @(24, "dog/dog084.png", "angeln", "dog"),
    @(126, "car/car097.png", "dienen", "car"),
    @(7, "dog/dog068.png", "herrschen", "dog"),
    @(95, "car/car066.png", "lassen", "car"),
    @(101, "dog/dog070.png", "wachsen", "dog"),
    @(19, "car/car087.png", "streifen", "car")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
